$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Controllers")
$ws.Range("E78").Value = "Ronald"
$ws.Range("E78").Style = $ws.Range("E86").Style
